$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for row 23 / new row 24 (columns A..AV), copied identically except the Date column (Y)
$values = @(565,474,450,522,515,520,474,570,490,450,571,480,485,505,545,480,618,490,474,480,619,550,599,495,45754,850,555,543.5,500,545,507,509,745,473,735,474,488,570,555,488,535,547,568,547,645,634,496,485)

# Column Y is column index 25 (date column)
$dateCol = 25

for ($c = 1; $c -le $values.Length; $c++) {
    $ws.Cells.Item(24, $c).Value = $values[$c - 1]
}

# The new row 24's Date cell (Y24) keeps the original date-only format that Y23 used to have
$ws.Cells.Item(24, $dateCol).NumberFormat = "YYYY-MM-DD"

# Y23 switches to the date+time number format
$ws.Cells.Item(23, $dateCol).NumberFormat = "YYYY-MM-DD HH:MM:SS"
